# NMCARS-PART-5215 edits
$d = $word.ActiveDocument

function Set-ParaXml($para, [string]$innerWordXml) {
    $pkg = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $innerWordXml + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
    $para.Range.InsertXML($pkg)
}

# --- 1) "(1) If a member of the armed forces..." paragraph ---
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "If a member of the armed forces*") {
        $inner = '<w:body><w:p w:rsidR="00C24F59" w:rsidRDefault="00C24F59" w:rsidP="00C24F59"><w:pPr><w:pStyle w:val="List2"/></w:pPr><w:r><w:rPr><w:szCs w:val="24"/></w:rPr><w:t>(1)</w:t></w:r><w:r><w:rPr><w:szCs w:val="24"/></w:rPr><w:tab/></w:r><w:r w:rsidRPr="00B46192"><w:t xml:space="preserve">If a member of the armed forces, is a flag or general officer; or </w:t></w:r></w:p></w:body>'
        Set-ParaXml $p $inner
        break
    }
}

# --- 2) "(2) If a civilian, is a member of the SES..." paragraph ---
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "If a civilian, is a member of the SES*") {
        $inner = '<w:body><w:p w:rsidR="00C24F59" w:rsidRDefault="00C24F59" w:rsidP="00C24F59"><w:pPr><w:pStyle w:val="List2"/></w:pPr><w:r><w:rPr><w:szCs w:val="24"/></w:rPr><w:t>(2)</w:t></w:r><w:r><w:rPr><w:szCs w:val="24"/></w:rPr><w:tab/></w:r><w:r w:rsidRPr="00B46192"><w:t>If a civilian, is a member of the SES (or in a com</w:t></w:r><w:r><w:t xml:space="preserve">parable or higher position under </w:t></w:r><w:r w:rsidRPr="00672AF7"><w:t>another schedule).</w:t></w:r></w:p></w:body>'
        Set-ParaXml $p $inner
        break
    }
}

# --- 3) "(4) Cost information..." paragraph: split "(4) " into "(4)" + " " runs ---
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Cost information*") {
        $inner = '<w:body><w:p w:rsidR="00C24F59" w:rsidRDefault="00C24F59" w:rsidP="00C24F59"><w:pPr><w:pStyle w:val="List2"/></w:pPr><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>(4)</w:t></w:r><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i/><w:color w:val="000000"/></w:rPr><w:t>Cost information</w:t></w:r><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>. The sharing of cost information with the technical evaluation team, and any limitations on the timing and extent of such sharing, should be addressed during the planning for the source selection. HCAs may establish specific procedural requirements for approving, documenting and/or varying from plans related to such sharing.</w:t></w:r></w:p></w:body>'
        Set-ParaXml $p $inner
        break
    }
}

# --- 4) "(2) After the contractor submits..." paragraph: pStyle Normalwline -> List2 ---
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*After the contractor submits*") {
        $p.Style = "List2"
        break
    }
}

# --- 5) "(c)(4)(A)(2) Senior procurement executive coordination" paragraph: add pStyle List1 ---
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*(c)(4)(A)(2)*") {
        $p.Style = "List1"
        break
    }
}

# --- 6) "(6) Submit a courtesy copy..." paragraph: pStyle List1 -> List2, split run "(6) Submit..." ---
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Submit a courtesy copy*") {
        $inner = '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p w:rsidR="00C24F59" w:rsidRDefault="00C24F59" w:rsidP="00C24F59"><w:pPr><w:pStyle w:val="List2"/></w:pPr><w:r><w:t>(6)</w:t></w:r><w:r><w:t xml:space="preserve"> Submit a courtesy copy of the quarterly report to DASN(P) by email </w:t></w:r><w:hyperlink r:id="rId15" w:history="1"><w:r w:rsidRPr="00AF1F53"><w:t>RDAJ&amp;As.fct@navy.mil</w:t></w:r></w:hyperlink><w:r><w:t xml:space="preserve"> with the subject &#8220;[Activity Name] DFARS 215.403-3 HCA Determination to Award to Offeror Failing to Comply with Requests for Data Other Than Certified Cost or Pricing Data.&#8221; Negative reports are required.</w:t></w:r></w:p></w:body></w:document>'
        $pkg = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' + $inner + '</pkg:xmlData></pkg:part></pkg:package>'
        $p.Range.InsertXML($pkg)
        $h = $p.Range.Hyperlinks(1)
        $h.Range.Style = "Hyperlink"
        break
    }
}

Write-Host "done"
